$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 34 (pushes existing rows 34..61 down to 35..62,
# and shifts the sheet dimension from A1:R61 to A1:R62).
$ws.Rows.Item(34).Insert()

# Populate the newly inserted row 34 with the new market observation.
$ws.Cells.Item(34, 1).Value = 10
$ws.Cells.Item(34, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(34, 3).Value = "La Araucanía"
$ws.Cells.Item(34, 4).Value = 44586
$ws.Cells.Item(34, 5).Value = 9
$ws.Cells.Item(34, 6).Value = 100112030
$ws.Cells.Item(34, 7).Value = "Poroto granado"
$ws.Cells.Item(34, 8).Value = "Sin especificar"
$ws.Cells.Item(34, 9).Value = "Primera"
$ws.Cells.Item(34, 10).Value = 55
$ws.Cells.Item(34, 11).Value = 28000
$ws.Cells.Item(34, 12).Value = 28000
$ws.Cells.Item(34, 13).Value = 28000
$ws.Cells.Item(34, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(34, 15).Value = "Región del Maule"
$ws.Cells.Item(34, 16).Value = 1120
$ws.Cells.Item(34, 17).Value = 25
$ws.Cells.Item(34, 18).Value = "Hortaliza"
